$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the existing "header/index" style (bold font, thin border all
# round, centered/top-aligned) from a cell that already carries it, far away
# from the table, so we can re-apply it after rebuilding the table. ---
$ws.Range("A2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)  # xlPasteFormats

# --- Clear only the table's original extent (leave the stash cell alone) ---
$ws.Range("A1:H5").Clear()

# --- Header row ---
$headers = @("peak_label","mz_mean","mz_width","rt","rt_min","rt_max","rt_unit","intensity_threshold","target_filename")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c-1]
}

# --- Data rows ---
# columns: peak_label, mz_mean, mz_width, rt, rt_min, rt_max, rt_unit, intensity_threshold, target_filename
$data = @(
    @("1", 151.0605, 5, $null, 304.2, 305.4, "s", 0, "v0.csv"),
    @("2", 216.0504, 5, $null, 238.8, 263.4, "s", 0, "v0.csv"),
    @("3", 115.0032, 5, $null, 207,   263.4, "s", 0, "v0.csv")
)

# peak_label column (A) holds text that looks numeric ("1","2","3") - force
# text storage so it round-trips as a shared string, not a number.
$ws.Range("A2:A4").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le 9; $c++) {
        $val = $row[$c-1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# --- Re-apply the bold/border/center style to the header row and the
# peak_label column, matching the original workbook's formatting ---
$ws.Range("Z100").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats

# --- cleanup stash cell ---
$ws.Range("Z100").Clear()

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
